$d = $word.ActiveDocument

# --- Change 1: subtitle -- split author line into two lines ----------------
# "Marlin Lee, University of Wisconsin Data Science Institute - July, 2022"
# becomes two runs (separated by a line break) -
#   "Marlin Lee, Steve Goldstein, Kyllan Wunder, Abe Megahed"
#   "University of Wisconsin Data Science Institute - July, 2022"
$subtitleOld = "Marlin Lee, University of Wisconsin Data Science Institute - July, 2022"

$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $subtitleOld) {
        $start = $para.Range.Start
        $end = $para.Range.End
        # exclude the trailing paragraph mark from the replaced range
        $target = $d.Range($start, $end - 1)

        $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t>Marlin Lee, Steve Goldstein, Kyllan Wunder, Abe Megahed</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t>University of Wisconsin Data Science Institute - July, 2022</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $target.InsertXML($xmlFrag)
        $found1 = $true
        break
    }
}
Write-Output "Subtitle split applied: $found1"

# --- Change 2: collapse a double space to a single space --------------------
$found2 = $d.Content.Find.Execute(
    "a lower wastewater sampling rate.   Wastewater facilities",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "a lower wastewater sampling rate. Wastewater facilities",
    2)
Write-Output "Double-space fix applied: $found2"

# --- Change 3: reword the closing sentence of the exponential-smoothing note
$oldTail = "we increase the slope between data points.  It would be possible to introduce a normalizing factor in order to make this determination more closely match the scale used for the original data."
$newTail = "we increase the slope between data points.    If we applied a normalizing factor to account for the higher slope due to the lower sampling rate, then the flagged regions would more closely match."

$found3 = $d.Content.Find.Execute(
    $oldTail,
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    $newTail,
    2)
Write-Output "Closing sentence reworded: $found3"
